$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 745.9375
$ws.Range("I4").Value = 454
$ws.Range("J4").Value = 1621.75
$ws.Range("K4").Value = 454
$ws.Range("L4").Value = 1621.75
$ws.Range("M4").Value = -340
$ws.Range("N4").Value = -1849.75
$ws.Range("H18").Value = 7966.3335
$ws.Range("I18").Value = 8759.6
$ws.Range("J18").Value = 4000
$ws.Range("K18").Value = 8759.6
$ws.Range("L18").Value = 4000
$ws.Range("M18").Value = -8475.6
$ws.Range("N18").Value = -4568
$ws.Range("H99").Value = 244.25
$ws.Range("I99").Value = 244.25
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 732.75
$ws.Range("L99").Value = 0
$ws.Range("M99").Value = 765.25
$ws.Range("N99").ClearContents()
$ws.Range("H106").Value = 2666.6667
$ws.Range("I106").Value = 2666.6667
$ws.Range("K106").Value = 2666.6667
$ws.Range("M106").Value = -2035.6667

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 4729.857
$ws.Range("I45").Value = 6035.8
$ws.Range("J45").Value = 1465
$ws.Range("K45").Value = 6035.8
$ws.Range("L45").Value = 1465
$ws.Range("M45").Value = -5658.8
$ws.Range("N45").Value = -2219
$ws.Range("H132").Value = 2323.5
$ws.Range("I132").Value = 2402.3333
$ws.Range("K132").Value = 7206.999899999999
$ws.Range("M132").Value = -4676.999899999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 323.66666
$ws.Range("I22").Value = 384.5
$ws.Range("J22").Value = 202
$ws.Range("K22").Value = 384.5
$ws.Range("L22").Value = 202
$ws.Range("M22").Value = -211.5
$ws.Range("N22").Value = -548
$ws.Range("H94").Value = 993.8
$ws.Range("J94").Value = 499.5
$ws.Range("L94").Value = 499.5
$ws.Range("N94").Value = -1401.5
$ws.Range("H105").Value = 2418.375
$ws.Range("I105").Value = 2392.4285
$ws.Range("J105").Value = 2600
$ws.Range("K105").Value = 2392.4285
$ws.Range("L105").Value = 2600
$ws.Range("M105").Value = -645.4285
$ws.Range("N105").Value = -6094
$ws.Range("H107").Value = 1542.6818
$ws.Range("I107").Value = 1702.3889
$ws.Range("J107").Value = 824
$ws.Range("K107").Value = 1702.3889
$ws.Range("L107").Value = 824
$ws.Range("M107").Value = 217.6111000000001
$ws.Range("N107").Value = -4664

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 6997.5
$ws.Range("I16").Value = 5996.6665
$ws.Range("K16").Value = 5996.6665
$ws.Range("M16").Value = -5709.6665
$ws.Range("H58").Value = 2245.5
$ws.Range("I58").Value = 2001.375
$ws.Range("K58").Value = 2001.375
$ws.Range("M58").Value = -1798.375
$ws.Range("H99").Value = 2666.4
$ws.Range("I99").Value = 2666.4
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 2666.4
$ws.Range("L99").Value = 0
$ws.Range("M99").Value = -1168.4
$ws.Range("N99").ClearContents()
$ws.Range("H113").Value = 6997.5
$ws.Range("I113").Value = 5996.6665
$ws.Range("K113").Value = 5996.6665
$ws.Range("M113").Value = -3826.6665
$ws.Range("H126").Value = 2666.4
$ws.Range("I126").Value = 2666.4
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 7999.200000000001
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -5529.200000000001
$ws.Range("N126").ClearContents()
$ws.Range("H132").Value = 4151.8184
$ws.Range("I132").Value = 4151.8184
$ws.Range("K132").Value = 12455.4552
$ws.Range("M132").Value = -9925.4552
$ws.Range("H134").Value = 3205.5
$ws.Range("I134").Value = 3028.5715
$ws.Range("K134").Value = 9085.7145
$ws.Range("M134").Value = -6550.7145
$ws.Range("H136").Value = 2245.5
$ws.Range("I136").Value = 2001.375
$ws.Range("K136").Value = 6004.125
$ws.Range("M136").Value = -3454.125
$ws.Range("H138").Value = 58178
$ws.Range("J138").Value = 58178
$ws.Range("L138").Value = 58178
$ws.Range("N138").Value = -68458

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H141").Value = 5721.1113
$ws.Range("I141").Value = 5721.1113
$ws.Range("J141").Value = 0
$ws.Range("K141").Value = 17163.3339
$ws.Range("L141").Value = 0
$ws.Range("M141").Value = -11983.3339
$ws.Range("N141").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H31").Value = 298.66666
$ws.Range("I31").Value = 298.66666
$ws.Range("K31").Value = 298.66666
$ws.Range("M31").Value = -6.666659999999979
$ws.Range("H37").Value = 298.66666
$ws.Range("I37").Value = 298.66666
$ws.Range("K37").Value = 298.66666
$ws.Range("M37").Value = -21.66665999999998
$ws.Range("H80").Value = 2995.3635
$ws.Range("I80").Value = 3750
$ws.Range("J80").Value = 2827.6667
$ws.Range("K80").Value = 3750
$ws.Range("L80").Value = 2827.6667
$ws.Range("M80").Value = -2752
$ws.Range("N80").Value = -4823.6667
$ws.Range("H83").Value = 2995.3635
$ws.Range("I83").Value = 3750
$ws.Range("J83").Value = 2827.6667
$ws.Range("K83").Value = 18750
$ws.Range("L83").Value = 14138.3335
$ws.Range("M83").Value = -13758
$ws.Range("N83").Value = -24122.3335
$ws.Range("H140").Value = 142712.25
$ws.Range("J140").Value = 142712.25
$ws.Range("L140").Value = 142712.25
$ws.Range("N140").Value = -153072.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 3252
$ws.Range("I40").Value = 2332.8
$ws.Range("J40").Value = 5550
$ws.Range("K40").Value = 2332.8
$ws.Range("L40").Value = 5550
$ws.Range("M40").Value = -2196.8
$ws.Range("N40").Value = -5822
$ws.Range("H82").Value = 1899.1177
$ws.Range("I82").Value = 2150.125
$ws.Range("J82").Value = 1676
$ws.Range("K82").Value = 2150.125
$ws.Range("L82").Value = 1676
$ws.Range("M82").Value = -1789.125
$ws.Range("N82").Value = -2398
$ws.Range("H85").Value = 1899.1177
$ws.Range("I85").Value = 2150.125
$ws.Range("J85").Value = 1676
$ws.Range("K85").Value = 2150.125
$ws.Range("L85").Value = 1676
$ws.Range("M85").Value = -902.125
$ws.Range("N85").Value = -4172
$ws.Range("H136").Value = 3323.25
$ws.Range("I136").Value = 1749.5
$ws.Range("J136").Value = 4897
$ws.Range("K136").Value = 5248.5
$ws.Range("L136").Value = 14691
$ws.Range("M136").Value = -2698.5
$ws.Range("N136").Value = -19791

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value = 98429
$ws.Range("J46").Value = 98429
$ws.Range("L46").Value = 98429
$ws.Range("N46").Value = -98891
$ws.Range("H132").Value = 2238.3635
$ws.Range("I132").Value = 2262.2
$ws.Range("K132").Value = 6786.599999999999
$ws.Range("M132").Value = -4256.599999999999
$ws.Range("H134").Value = 98429
$ws.Range("J134").Value = 98429
$ws.Range("L134").Value = 295287
$ws.Range("N134").Value = -300357
$ws.Range("H136").Value = 3148.0952
$ws.Range("I136").Value = 3345.111
$ws.Range("J136").Value = 1966
$ws.Range("K136").Value = 10035.333
$ws.Range("L136").Value = 5898
$ws.Range("M136").Value = -7485.332999999999
$ws.Range("N136").Value = -10998
